# Add files via upload
# The author re-uploaded the workbook with a new "Project Status" column
# inserted before the existing "Actual Cost to Date (Mil)" column. Net
# effect on the data: column Q (header + 86 rows) now holds what used to
# be in column R ("Project Status": In Progress/Done), and column R now
# holds what used to be in column Q (the cost figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the contents of columns Q and R (header row 1 + data rows 2:86) ---
$qVals = $ws.Range("Q1:Q86").Value()
$rVals = $ws.Range("R1:R86").Value()

$ws.Range("Q1:Q86").Value = $rVals
$ws.Range("R1:R86").Value = $qVals

# --- Column widths for the now-wider "Project Status" / cost columns ---
$ws.Columns("Q").ColumnWidth = 11.66666667
$ws.Columns("R").ColumnWidth = 20.83333333

# --- The filter-database defined name now spans through column R ---
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$R`$86"

# --- Selection moved to T4 ---
[void]$ws.Range("T4").Select()
